# Weekly update for the "Hortaliza, Vega Central Mapocho de Santiago - Cilantro" sheet.
#
# The published series picked up a new week's worth of readings (two new
# observations dated 2021-11-08 / serial 44508), which get inserted at the
# top of the "Primera" quality price block (rows 306-307). Every record that
# used to live at row N >= 306 shifts down two rows to make room, so the
# sheet grows from 381 to 383 data+header rows overall.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 306:381 (and everything after) down by two rows, leaving two
# blank rows at 306:307 ready for the new observations.
$ws.Range("A306:R307").Insert()

# New "caja 36 atados" reading for 2021-11-08.
$row306 = @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44508, 13, 100112040, "Cilantro", "Sin especificar", "Primera", 34, 4000, 4500, 4250, "`$/caja 36 atados", "Región Metropolitana", 118, 36, "Hortaliza")

# New "docena de atados" reading for 2021-11-08.
$row307 = @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44508, 13, 100112040, "Cilantro", "Sin especificar", "Primera", 106, 8000, 10000, 9000, "`$/docena de atados", "Región Metropolitana", 3000, 3, "Hortaliza")

for ($i = 0; $i -lt $row306.Length; $i++) {
    $ws.Cells.Item(306, $i + 1).Value = $row306[$i]
    $ws.Cells.Item(307, $i + 1).Value = $row307[$i]
}
